# excel_writer: also include totals for the balance columns
#
# The "Gesamtergebnis" (totals) sheet previously left the balance columns
# (Startguthaben / Endsaldo) as "N/A" on the Total row. Now they are
# included in the totals as numeric 0 values, like all the other columns.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Gesamtergebnis")

# Totals row (row 3): balance columns C (Startguthaben) and D (Endsaldo)
# now report a numeric total instead of "N/A".
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 0

# Row height for the totals row shrinks slightly now that the cells hold
# plain numbers rather than text.
$ws3.Rows.Item(3).RowHeight = 13.8

# The totals sheet becomes the active sheet/selection on save.
$ws3.Activate() | Out-Null
$ws3.Range("C6").Select() | Out-Null
